# ZRADTriBOM.xlsx - Release Rev B
# - Replace the D1 RGB LED (QLSP14RGB_B / Quelighting) with a more common
#   part (BL-HJXGXBX32M-D / American Bright Optoelectronics Corp), updating
#   its cost and adding a caution note about alternates/pinouts.
# - Clarify the optional LED light-pipe comment on the D2 row with the
#   recommended mounting-hole size.
# - Log the change in the Introduction sheet's revision history table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# ZRADTriBOM sheet: update the D1 RGB LED line item (row 16)
# ---------------------------------------------------------------------
$bom = $wb.Worksheets.Item("ZRADTriBOM")

$bom.Range("B16").Value = "BL-HJXGXBX32M-D"
$bom.Range("C16").Value = "American Bright Optoelectronics Corp"
$bom.Range("D16").Value = "RGB LED 3.2x2.8mm "
$bom.Range("F16").Value = "BL-HJXGXBX32M-DCT-ND"
$bom.Range("H16").Value = 0.17
$bom.Range("N16").Value = "CAUTION! Carefully choose alternates as there are many flavors of pinouts!"

# Clarify the optional LED light pipe comment (row 40, D2 section)
$bom.Range("M40").Value = "Optional - makes the LED much more visible - 5/32 or 4mm hole"

# Move the active selection to reflect where the edit was made
$bom.Activate()
$bom.Range("C16").Select()

# ---------------------------------------------------------------------
# Introduction sheet: log the change in the revision history table
# ---------------------------------------------------------------------
$intro = $wb.Worksheets.Item("Introduction")

# Copy the formatting of the prior history row down to the new row so the
# date keeps the same date style/number format.
$intro.Range("A17").Copy($intro.Range("A18"))

$intro.Range("A18").Value = 45839
$intro.Range("B18").Value = "Changed RGB LED to 3.2x2.8mm size which is more common"

$intro.Activate()
$intro.Range("B18").Select()

# Leave ZRADTriBOM as the active/displayed tab, matching the released file.
$bom.Activate()

$wb.Application.CalculateFull()
